$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Renumber column E for rows 13-17 (the "No." column), shifting values
# down by one starting at row 13 so the sequence reads 1..16 without
# skipping "12" (it previously jumped from 11 straight to 13).
# These "No." values are stored as text throughout the sheet, so a
# leading apostrophe is used to enter them as text rather than numbers.
$ws.Range("E13").Value = "'12"
$ws.Range("E14").Value = "'13"
$ws.Range("E15").Value = "'14"
$ws.Range("E16").Value = "'15"
$ws.Range("E17").Value = "'16"
